$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# --- New data in columns B:E, rows 56-61 (feature selection / PCA results) ---
# Set the raw values first.
$ws.Range("B56").Value = "train_mae"

$ws.Range("C56").Value = "SVR"
$ws.Range("D56").Value = 0.0053
$ws.Range("E56").Value = 0.1

$ws.Range("C57").Value = "lightgbm"
$ws.Range("D57").Value = 0.0562
$ws.Range("E57").Value = 0.05

$ws.Range("C58").Value = "xgboost"
$ws.Range("D58").Value = 0.055
$ws.Range("E58").Value = 0.1

$ws.Range("C59").Value = "ridge"
$ws.Range("D59").Value = 0
$ws.Range("E59").Value = 0.25

$ws.Range("C60").Value = "rf"
$ws.Range("D60").Value = 0.0174
$ws.Range("E60").Value = 0.1

$ws.Range("C61").Value = "gbr"
$ws.Range("D61").Value = 0.0587
$ws.Range("E61").Value = 0.05

# Row 62 - sum formula mirroring column O's SUM(O56:O61)
$ws.Range("E62").Formula = "=SUM(E56:E61)"

# --- Copy the existing formatting from the matching M/N columns (rows 56-61)
# so the new C/D columns reuse the exact same cell styles (Consolas left/center
# aligned text, Courier New numbers) instead of creating new ones. ---
$ws.Range("M56:M61").Copy()
$ws.Range("C56:C61").PasteSpecial(-4122)

$ws.Range("N56:N61").Copy()
$ws.Range("D56:D61").PasteSpecial(-4122)

# --- Update sheet view: scrolled position and active selection cell ---
$excel.ActiveWindow.ScrollRow = 46
$ws.Range("E56").Select()
